$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.952.64'
$ws.Range("E2").Value = '  +0.81%  '
$ws.Range("D3").Value = '3.395.98'
$ws.Range("E3").Value = '  +0.53%  '
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = "'580.46"
$ws.Range("E5").Value = '  +0.64%  '
$ws.Range("D6").Value = "'138.47"
$ws.Range("E6").Value = '  +1.84%  '
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = '  -0.13%  '
$ws.Range("D8").Value = '3.394.75'
$ws.Range("E8").Value = '  +0.57%  '
$ws.Range("E9").Value = '  -0.68%  '
$ws.Range("D10").Value = "'7.54"
$ws.Range("E10").Value = '  -0.29%  '
$ws.Range("E11").Value = '  +2.83%  '
$ws.Range("E12").Value = '  +1.06%  '
$ws.Range("D13").Value = '3.975.60'
$ws.Range("E13").Value = '  +0.51%  '
$ws.Range("E14").Value = '  +1.66%  '
$ws.Range("E15").Value = '  +1.84%  '
$ws.Range("D16").Value = '3.396.83'
$ws.Range("E16").Value = '  +0.52%  '
$ws.Range("E17").Value = '  +0.80%  '
$ws.Range("D18").Value = '62.045.85'
$ws.Range("E18").Value = '  +0.74%  '
$ws.Range("E19").Value = '  +1.27%  '
$ws.Range("E20").Value = '  +1.52%  '
$ws.Range("E21").Value = '  +0.22%  '
$ws.Range("D22").Value = "'390.47"
$ws.Range("E22").Value = '  +1.36%  '
$ws.Range("E23").Value = '  -0.40%  '
$ws.Range("E24").Value = '  +8.61%  '
$ws.Range("D25").Value = '3.537.80'
$ws.Range("E25").Value = '  +0.70%  '
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = '  -0.16%  '
$ws.Range("D27").Value = "'71.62"
$ws.Range("E27").Value = '  +0.87%  '
$ws.Range("D28").Value = "'1.68"
$ws.Range("E28").Value = '  -1.61%  '
$ws.Range("D29").Value = "'7.66"
$ws.Range("E29").Value = '  -2.39%  '
$ws.Range("E30").Value = '  -0.02%  '
$ws.Range("D31").Value = "'0.161"
$ws.Range("E31").Value = '  +2.93%  '
$ws.Range("D32").Value = "'8.25"
$ws.Range("E32").Value = '  +1.08%  '
$ws.Range("D33").Value = "'2.18"
$ws.Range("E33").Value = '  +0.29%  '
$ws.Range("E34").Value = '  +0.00%  '
$ws.Range("D35").Value = "'23.53"
$ws.Range("D36").Value = '3.428.15'
$ws.Range("E36").Value = '  +0.48%  '
$ws.Range("D37").Value = "'5.39"
$ws.Range("E37").Value = '  -3.40%  '
$ws.Range("E38").Value = '  +2.45%  '
$ws.Range("D39").Value = "'6.89"
$ws.Range("E39").Value = '  -1.24%  '
$ws.Range("D40").Value = "'165.20"
$ws.Range("E40").Value = '  +1.53%  '
$ws.Range("D41").Value = "'0.0789"
$ws.Range("E41").Value = '  +0.24%  '
$ws.Range("D42").Value = "'1.77"
$ws.Range("E42").Value = '  +9.34%  '
$ws.Range("E43").Value = '  +2.08%  '
$ws.Range("E44").Value = '  +3.82%  '
$ws.Range("E45").Value = '  -0.03%  '
$ws.Range("D46").Value = "'25.31"
$ws.Range("E46").Value = '  +7.73%  '
$ws.Range("D47").Value = "'4.44"
$ws.Range("E47").Value = '  -0.07%  '
$ws.Range("D48").Value = "'41.40"
$ws.Range("E48").Value = '  -0.67%  '
$ws.Range("D49").Value = "'6.90"
$ws.Range("E49").Value = '  -0.26%  '
$ws.Range("D50").Value = "'23.17"
$ws.Range("E50").Value = '  -0.17%  '
$ws.Range("D51").Value = '2.351.03'
$ws.Range("E51").Value = '  +6.84%  '
